# Applies the cryptos list refresh for Thu Jun  6 14:45:07 UTC 2024 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to remain plain text,
# even when the text looks like a number (e.g. "36.83"). The sheet stores
# these price/volume columns as text, so we briefly flip the cell to the
# Text number format before assigning, then restore the default "Normal"
# style so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "71.268.01"
$ws.Range("E2").Value = "  +0.48%  "

Set-TextValue $ws.Range("D3") "3.841.52"
$ws.Range("E3").Value = "  +1.05%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue $ws.Range("D5") "714.21"
$ws.Range("E5").Value = "  +1.72%  "

Set-TextValue $ws.Range("D6") "173.01"
$ws.Range("E6").Value = "  +0.17%  "

Set-TextValue $ws.Range("D7") "3.840.40"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("E13").Value = "  +0.14%  "

Set-TextValue $ws.Range("D14") "36.83"
$ws.Range("E14").Value = "  +2.09%  "

Set-TextValue $ws.Range("D15") "4.490.61"
$ws.Range("E15").Value = "  +1.09%  "

Set-TextValue $ws.Range("D16") "3.839.34"
$ws.Range("E16").Value = "  +0.67%  "

Set-TextValue $ws.Range("D17") "71.260.63"
$ws.Range("E17").Value = "  +0.63%  "

Set-TextValue $ws.Range("D18") "7.24"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("E19").Value = "  +0.47%  "

Set-TextValue $ws.Range("D20") "17.43"
$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("E21").Value = "  -1.95%  "

Set-TextValue $ws.Range("D22") "495.72"
$ws.Range("E22").Value = "  +2.99%  "

$ws.Range("E23").Value = "  +2.32%  "

Set-TextValue $ws.Range("D24") "85.30"
$ws.Range("E24").Value = "  +1.37%  "

$ws.Range("E25").Value = "  +1.79%  "

Set-TextValue $ws.Range("D26") "10.68"
$ws.Range("E26").Value = "  +1.59%  "

$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D28") "2.10"
$ws.Range("E28").Value = "  -2.74%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D29") "3.19"
$ws.Range("E29").Value = "  +2.45%  "

$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("E32").Value = "  -2.33%  "

$ws.Range("E33").Value = "  +0.04%  "

Set-TextValue $ws.Range("D34") "0.180"
$ws.Range("E34").Value = "  -4.79%  "

Set-TextValue $ws.Range("D35") "9.22"
$ws.Range("E35").Value = "  -0.46%  "

Set-TextValue $ws.Range("D36") "3.806.93"
$ws.Range("E36").Value = "  +1.51%  "

Set-TextValue $ws.Range("D37") "0.996"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("E39").Value = "  +5.70%  "

Set-TextValue $ws.Range("D40") "6.02"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  -1.27%  "

$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("E44").Value = "  +0.19%  "

Set-TextValue $ws.Range("D45") "0.000321"
$ws.Range("E45").Value = "  +1.47%  "

Set-TextValue $ws.Range("D46") "163.96"
$ws.Range("E46").Value = "  -0.10%  "

Set-TextValue $ws.Range("D47") "48.91"
$ws.Range("E47").Value = "  -0.03%  "

Set-TextValue $ws.Range("D48") "424.81"
$ws.Range("E48").Value = "  +3.66%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("E51").Value = "  -1.11%  "
